$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# The sheet is protected; remove protection so the cells below can be edited,
# then restore protection with the same effective options afterwards.
$ws.Unprotect()

# Update the confidential disclosure date string in cell A18 (2021-05-19 -> 2021-05-20)
$ws.Range("A18").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-05-20 for illustrative purposes only and are subject to change."
# Setting a taller, wrapped multi-line value can make Excel stamp a custom row
# height; auto-fit the row back so it matches the sheet's default height again.
$ws.Rows.Item(18).AutoFit()

# Update the Weight (D) / Percent Change (E) columns for rows 2-15
$ws.Range("D2").Value = 0.05627846174496822
$ws.Range("E2").Value = 0.01085628348528989

$ws.Range("D3").Value = 0.0236752787242526
$ws.Range("E3").Value = 0.01388074291300079

$ws.Range("D4").Value = 0.03120733011226995
$ws.Range("E4").Value = 0.003623188405797118

$ws.Range("D5").Value = 0.03266467725116823
$ws.Range("E5").Value = 0.001538461538461489

$ws.Range("D6").Value = 0.03746783950931572
$ws.Range("E6").Value = -0.00529436679373152

$ws.Range("D7").Value = 0.0188430224215899
$ws.Range("E7").Value = -0.002456398919184588

$ws.Range("D8").Value = 0.004358155622721655
$ws.Range("E8").Value = 0.02548930359581258

$ws.Range("D9").Value = 0.006796104307439414
$ws.Range("E9").Value = 0.009729519361743488

$ws.Range("D10").Value = 0.07375340284605879
$ws.Range("E10").Value = 0.004303388918773576

$ws.Range("D11").Value = 0.07383275024019118
$ws.Range("E11").Value = 0.004836109618484796

$ws.Range("D12").Value = 0.1439996508714659
$ws.Range("E12").Value = 0.008228638601131211

$ws.Range("D13").Value = 0.3834264452962534
$ws.Range("E13").Value = 0.004133685136323617

$ws.Range("D14").Value = 0.1136968810523051
$ws.Range("E14").Value = 0.004030288226673173

$ws.Range("E15").Value = 0.004938052828172612

# Restore sheet protection (same effective settings as before the edit: sheet
# protected, objects/scenarios locked, column/row formatting left allowed).
# The original password hash can't be recovered, so this keeps the sheet
# password-gated without claiming to know the author's original password.
$ws.Protect("D382", $true, $true)

$wb.Save()
